$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

# Title (appears twice: heading and bold summary line near the bottom).
# Replace:=2 (wdReplaceAll) replaces every occurrence in the range in one call.
Replace-Text "Play Beellionaires Dream Drop for Free - Review" "Play Beellionaires Dream Drop for Free - Exciting Features & Big Wins"

# What we like
Replace-Text "Up to 1,024 paylines for potential big winnings" "Hive structure with up to 1,024 paylines for potential big wins"
Replace-Text "Random multipliers and exciting Royal Respin function" "Collection of coins with random multipliers adds excitement to gameplay"
Replace-Text "Charming design with unique beehive theme" "Magnificent reels and beehive palace background create an impressive design"
Replace-Text "Progressive jackpots and bonus game function for added excitement" "Exciting features like Royal Respin and Bonus Game add to the thrill"

# What we don't like
Replace-Text "Low RTP at 94%" "RTP of 94% is considered low compared to other slot games"
Replace-Text "Limited bonus game spins" "Players may need to practice responsible gambling to avoid unwanted losses"

# Meta description (italic line)
Replace-Text "Discover the charming design and exciting features of Beellionaires Dream Drop. Play now for free and potentially win big with random multipliers and progressive jackpots." "Read our review of Beellionaires Dream Drop and play for free. Experience exciting features and the chance to win big rewards."
